$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# New rows of data (TimeStamp + 14 numeric columns B..O) to append after row 52
$newRows = @(
    @(45701.931446759256, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116),
    @(45701.932812500003, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116),
    @(45702.475960648146, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116),
    @(45702.479409722226, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116),
    @(45702.480138888888, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116),
    @(45702.511018518519, 8, 6, 216, 389, 365, 388, 2681, 388, 1216, 119, 304, 30, 3111, 4116)
)

$startRow = 53
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}

$wb.Save()
